$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 468 ("この地上に生き甲斐あり...") which shifts all
# subsequent rows up by one.
$ws.Rows.Item(468).Delete()
